$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.081.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.780.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3805"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3399"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.426"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.777.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.066"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06638"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.36%  "

$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.510"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.082.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.371"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.485"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.440"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "155.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.979.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.976"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.993"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08661"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.628"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.60%  "

$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6809"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.44%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.365"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06272"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.528"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.40%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6416"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.851"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.115"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07084"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.44%  "
